# Weekly refresh of the Ciruela (plum) price sheet: each reporting date's
# pair of rows (Primera/Segunda quality) gets re-dated and its volume /
# price / packaging figures updated to the new week's values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44335
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 10000
$ws.Range("O2").Value = 11000
$ws.Range("P2").Value = 10500
$ws.Range("Q2").Value = '$/bandeja 18 kilos granel'
$ws.Range("S2").Value = 583
$ws.Range("T2").Value = 18

# Row 3
$ws.Range("D3").Value = 44335
$ws.Range("M3").Value = 50
$ws.Range("N3").Value = 9000
$ws.Range("O3").Value = 9000
$ws.Range("P3").Value = 9000
$ws.Range("Q3").Value = '$/bandeja 18 kilos granel'
$ws.Range("T3").Value = 18

# Row 4
$ws.Range("D4").Value = 44307
$ws.Range("K4").Value = 'Angeleno'
$ws.Range("N4").Value = 9000
$ws.Range("O4").Value = 10000
$ws.Range("P4").Value = 9500
$ws.Range("Q4").Value = '$/bandeja 18 kilos granel'
$ws.Range("S4").Value = 528
$ws.Range("T4").Value = 18

# Row 5
$ws.Range("D5").Value = 44307
$ws.Range("K5").Value = 'Angeleno'
$ws.Range("N5").Value = 8000
$ws.Range("O5").Value = 8000
$ws.Range("P5").Value = 8000
$ws.Range("Q5").Value = '$/bandeja 18 kilos granel'
$ws.Range("S5").Value = 444
$ws.Range("T5").Value = 18

# Row 6
$ws.Range("D6").Value = 44223
$ws.Range("M6").Value = 100
$ws.Range("N6").Value = 10000
$ws.Range("O6").Value = 11000
$ws.Range("P6").Value = 10500
$ws.Range("Q6").Value = '$/caja 16 kilos granel'
$ws.Range("S6").Value = 656
$ws.Range("T6").Value = 16

# Row 7
$ws.Range("D7").Value = 44223
$ws.Range("M7").Value = 50
$ws.Range("N7").Value = 9000
$ws.Range("O7").Value = 9000
$ws.Range("P7").Value = 9000
$ws.Range("Q7").Value = '$/caja 16 kilos granel'
$ws.Range("S7").Value = 562
$ws.Range("T7").Value = 16

# Row 8
$ws.Range("D8").Value = 44202
$ws.Range("K8").Value = 'Black Amber'
$ws.Range("N8").Value = 14000
$ws.Range("O8").Value = 15000
$ws.Range("P8").Value = 14500
$ws.Range("Q8").Value = '$/caja 18 kilos granel'
$ws.Range("S8").Value = 806

# Row 9
$ws.Range("D9").Value = 44202
$ws.Range("K9").Value = 'Black Amber'
$ws.Range("N9").Value = 12000
$ws.Range("O9").Value = 12000
$ws.Range("P9").Value = 12000
$ws.Range("Q9").Value = '$/caja 18 kilos granel'
$ws.Range("S9").Value = 667

# Row 10
$ws.Range("D10").Value = 44343
$ws.Range("K10").Value = 'Angeleno'
$ws.Range("M10").Value = 200
$ws.Range("Q10").Value = '$/bandeja 18 kilos granel'
$ws.Range("S10").Value = 583
$ws.Range("T10").Value = 18

# Row 11
$ws.Range("D11").Value = 44343
$ws.Range("K11").Value = 'Angeleno'
$ws.Range("M11").Value = 100
$ws.Range("Q11").Value = '$/bandeja 18 kilos granel'
$ws.Range("S11").Value = 500
$ws.Range("T11").Value = 18

# Row 12
$ws.Range("D12").Value = 44189
$ws.Range("K12").Value = 'Red Beaut'
$ws.Range("N12").Value = 12000
$ws.Range("O12").Value = 13000
$ws.Range("P12").Value = 12500
$ws.Range("Q12").Value = '$/caja 15 kilos granel'
$ws.Range("S12").Value = 833
$ws.Range("T12").Value = 15

# Row 13
$ws.Range("D13").Value = 44189
$ws.Range("K13").Value = 'Red Beaut'
$ws.Range("N13").Value = 10000
$ws.Range("O13").Value = 10000
$ws.Range("P13").Value = 10000
$ws.Range("Q13").Value = '$/caja 15 kilos granel'
$ws.Range("S13").Value = 667
$ws.Range("T13").Value = 15

# Row 14
$ws.Range("D14").Value = 44215
$ws.Range("N14").Value = 10000
$ws.Range("O14").Value = 11000
$ws.Range("P14").Value = 10500
$ws.Range("Q14").Value = '$/caja 16 kilos granel'
$ws.Range("S14").Value = 656
$ws.Range("T14").Value = 16

# Row 15
$ws.Range("D15").Value = 44215
$ws.Range("N15").Value = 8000
$ws.Range("O15").Value = 8000
$ws.Range("P15").Value = 8000
$ws.Range("Q15").Value = '$/caja 16 kilos granel'
$ws.Range("S15").Value = 500
$ws.Range("T15").Value = 16

# Row 16
$ws.Range("D16").Value = 44236
$ws.Range("K16").Value = 'Lemon'
$ws.Range("N16").Value = 14000
$ws.Range("O16").Value = 15000
$ws.Range("P16").Value = 14500
$ws.Range("Q16").Value = '$/caja 16 kilos granel'
$ws.Range("S16").Value = 906
$ws.Range("T16").Value = 16

# Row 17
$ws.Range("D17").Value = 44236
$ws.Range("K17").Value = 'Lemon'
$ws.Range("N17").Value = 12000
$ws.Range("O17").Value = 12000
$ws.Range("P17").Value = 12000
$ws.Range("Q17").Value = '$/caja 16 kilos granel'
$ws.Range("S17").Value = 750
$ws.Range("T17").Value = 16

# Row 18
$ws.Range("D18").Value = 44299
$ws.Range("M18").Value = 200
$ws.Range("N18").Value = 12000
$ws.Range("O18").Value = 13000
$ws.Range("P18").Value = 12500
$ws.Range("S18").Value = 694

# Row 19
$ws.Range("D19").Value = 44299
$ws.Range("M19").Value = 100
$ws.Range("N19").Value = 11000
$ws.Range("O19").Value = 11000
$ws.Range("P19").Value = 11000
$ws.Range("S19").Value = 611

# Row 22
$ws.Range("D22").Value = 44323
$ws.Range("K22").Value = 'Angeleno'
$ws.Range("N22").Value = 11000
$ws.Range("O22").Value = 12000
$ws.Range("P22").Value = 11500
$ws.Range("Q22").Value = '$/bandeja 18 kilos granel'
$ws.Range("S22").Value = 639
$ws.Range("T22").Value = 18

# Row 23
$ws.Range("D23").Value = 44323
$ws.Range("K23").Value = 'Angeleno'
$ws.Range("N23").Value = 9000
$ws.Range("O23").Value = 9000
$ws.Range("P23").Value = 9000
$ws.Range("Q23").Value = '$/bandeja 18 kilos granel'
$ws.Range("S23").Value = 500
$ws.Range("T23").Value = 18

# Row 24
$ws.Range("D24").Value = 44266
$ws.Range("K24").Value = 'Black Amber'

# Row 25
$ws.Range("D25").Value = 44266
$ws.Range("K25").Value = 'Black Amber'

# Row 26
$ws.Range("D26").Value = 44218

# Row 27
$ws.Range("D27").Value = 44218
$ws.Range("N27").Value = 9000
$ws.Range("O27").Value = 9000
$ws.Range("P27").Value = 9000
$ws.Range("S27").Value = 562

# Row 28
$ws.Range("D28").Value = 44328
$ws.Range("K28").Value = 'Angeleno'
$ws.Range("M28").Value = 100
$ws.Range("N28").Value = 9000
$ws.Range("O28").Value = 10000
$ws.Range("P28").Value = 9500
$ws.Range("Q28").Value = '$/bandeja 18 kilos granel'
$ws.Range("S28").Value = 528
$ws.Range("T28").Value = 18

# Row 29
$ws.Range("D29").Value = 44328
$ws.Range("K29").Value = 'Angeleno'
$ws.Range("M29").Value = 50
$ws.Range("N29").Value = 8000
$ws.Range("O29").Value = 8000
$ws.Range("P29").Value = 8000
$ws.Range("Q29").Value = '$/bandeja 18 kilos granel'
$ws.Range("S29").Value = 444
$ws.Range("T29").Value = 18

# Row 30
$ws.Range("D30").Value = 44285
$ws.Range("Q30").Value = '$/caja 18 kilos granel'

# Row 31
$ws.Range("D31").Value = 44285
$ws.Range("Q31").Value = '$/caja 18 kilos granel'

# Row 32
$ws.Range("D32").Value = 44251
$ws.Range("N32").Value = 9000
$ws.Range("O32").Value = 10000
$ws.Range("P32").Value = 9500
$ws.Range("Q32").Value = '$/caja 16 kilos granel'
$ws.Range("S32").Value = 594
$ws.Range("T32").Value = 16

# Row 33
$ws.Range("D33").Value = 44251
$ws.Range("N33").Value = 8000
$ws.Range("O33").Value = 8000
$ws.Range("P33").Value = 8000
$ws.Range("Q33").Value = '$/caja 16 kilos granel'
$ws.Range("S33").Value = 500
$ws.Range("T33").Value = 16

# Row 34
$ws.Range("D34").Value = 44279
$ws.Range("K34").Value = 'Black Amber'
$ws.Range("N34").Value = 9000
$ws.Range("O34").Value = 10000
$ws.Range("P34").Value = 9500
$ws.Range("Q34").Value = '$/bandeja 18 kilos granel'
$ws.Range("S34").Value = 528
$ws.Range("T34").Value = 18

# Row 35
$ws.Range("D35").Value = 44279
$ws.Range("K35").Value = 'Black Amber'
$ws.Range("N35").Value = 8000
$ws.Range("O35").Value = 8000
$ws.Range("P35").Value = 8000
$ws.Range("Q35").Value = '$/bandeja 18 kilos granel'
$ws.Range("S35").Value = 444
$ws.Range("T35").Value = 18

# Row 36
$ws.Range("D36").Value = 44257
$ws.Range("K36").Value = 'Black Amber'
$ws.Range("M36").Value = 200
$ws.Range("Q36").Value = '$/caja 15 kilos granel'
$ws.Range("S36").Value = 700
$ws.Range("T36").Value = 15

# Row 37
$ws.Range("D37").Value = 44257
$ws.Range("K37").Value = 'Black Amber'
$ws.Range("M37").Value = 100
$ws.Range("Q37").Value = '$/caja 15 kilos granel'
$ws.Range("S37").Value = 600
$ws.Range("T37").Value = 15
